# Estadisticos Segundo Parcial 26 Mayo
$wb = $excel.ActiveWorkbook

# --- Sheet "Calificaciones": update 2P grades for "Conservación de la energía..." (col K)
# and Final grades for the same subject (col AA) ---
$wsCal = $wb.Worksheets.Item("Calificaciones")
$wsCal.Range("K4").Value = 6
$wsCal.Range("K5").Value = 6
$wsCal.Range("K6").Value = 5
$wsCal.Range("K7").Value = 8
$wsCal.Range("K8").Value = 8
$wsCal.Range("AA8").Value = 9
$wsCal.Range("K9").Value = 6
$wsCal.Range("K10").Value = 9
$wsCal.Range("K11").Value = 6
$wsCal.Range("AA11").Value = 8
$wsCal.Range("K12").Value = 6
$wsCal.Range("K13").Value = 6
$wsCal.Range("K14").Value = 6

# --- Sheet "Asistencias": update F2 and F3 attendance percentages for
# "Conservación de la energía..." (cols K and S mirror each other) ---
$wsAsi = $wb.Worksheets.Item("Asistencias")
$wsAsi.Range("K5").Value = 97.7
$wsAsi.Range("S5").Value = 97.7
$wsAsi.Range("K6").Value = 67.40000000000001
$wsAsi.Range("S6").Value = 67.40000000000001
$wsAsi.Range("K9").Value = 72.09999999999999
$wsAsi.Range("S9").Value = 72.09999999999999
$wsAsi.Range("K11").Value = 90.7
$wsAsi.Range("S11").Value = 90.7
$wsAsi.Range("K12").Value = 79.09999999999999
$wsAsi.Range("S12").Value = 79.09999999999999
$wsAsi.Range("K13").Value = 95.3
$wsAsi.Range("S13").Value = 95.3
$wsAsi.Range("K14").Value = 79.09999999999999
$wsAsi.Range("S14").Value = 79.09999999999999

# --- Sheet "Totales": update average grade for
# "Conservación de la energía..." (H2) ---
$wsTot = $wb.Worksheets.Item("Totales")
$wsTot.Range("H2").Value = 6.5
